$wb = $excel.ActiveWorkbook

# ---- 1. Replace TC001-A sheet (position 2) with new minimal content ----
$old = $wb.Worksheets.Item("TC001-A")
$old.Delete()
$tc001 = $wb.Worksheets.Add()
$tc001.Name = "TC001-A"

$tc001.Range("A1").Value = "Click"
$tc001.Range("B1").Value = "Id"
$tc001.Range("C1").Value = "et_auth_code"
$tc001.Range("D1").Value = "N/A"
$tc001.Range("E1").Value = "Click in the Continue button"

$tc001.Range("A2").Value = "SendText"
$tc001.Range("B2").Value = "Id"
$tc001.Range("C2").Value = "et_auth_code"
$tc001.Range("D2").Value = "Test"
$tc001.Range("E2").Value = "Send keys to auth code text field"

$tc001.Range("A1:E2").Select()

# ---- 2. Rewrite TC002-i (sheet3) with expanded content ----
$tc002 = $wb.Worksheets.Item("TC002-i")

$tc002.Range("A1").Value = "Click"
$tc002.Range("B1").Value = "Id"
$tc002.Range("C1").Value = "Allow"
$tc002.Range("D1").Value = "N/A"
$tc002.Range("E1").Value = "Click in the Continue button"

$tc002.Range("A2").Value = "Click"
$tc002.Range("B2").Value = "name"
$tc002.Range("C2").Value = "//XCUIElementTypeTextField"
$tc002.Range("D2").Value = "N/A"
$tc002.Range("E2").Value = "Tap on auth code field"

$tc002.Range("A3").Value = "SendText"
$tc002.Range("B3").Value = "name"
$tc002.Range("C3").Value = "//XCUIElementTypeTextField"
$tc002.Range("D3").Value = "test"
$tc002.Range("E3").Value = "Type authentication code into the text field"

$tc002.Range("A4").Value = "Click"
$tc002.Range("B4").Value = "name"
$tc002.Range("C4").Value = "Continue"
$tc002.Range("D4").Value = "N/A"
$tc002.Range("E4").Value = "Click in the continue button"

$tc002.Range("A5").Value = "Click"
$tc002.Range("B5").Value = "Xpath"
$tc002.Range("C5").Value = '//XCUIElementTypeTextField[contains(@value, "Username")]'
$tc002.Range("D5").Value = "N/A"
$tc002.Range("E5").Value = "Tap on username text field"

$tc002.Range("A6").Value = "SendText"
$tc002.Range("B6").Value = "Xpath"
$tc002.Range("C6").Value = '//XCUIElementTypeTextField[contains(@value, "Username")]'
$tc002.Range("D6").Value = "test.admin@ur.com"
$tc002.Range("E6").Value = "Type username into the text field"

$tc002.Range("A7").Value = "Click"
$tc002.Range("B7").Value = "Xpath"
$tc002.Range("C7").Value = "//XCUIElementTypeSecureTextField[contains(@value, Password)]"
$tc002.Range("D7").Value = "N/A"
$tc002.Range("E7").Value = "Tap on password text field"

$tc002.Range("A8").Value = "SendText"
$tc002.Range("B8").Value = "Xpath"
$tc002.Range("C8").Value = "//XCUIElementTypeSecureTextField[contains(@value, Password)]"
$tc002.Range("D8").Value = "test"
$tc002.Range("E8").Value = "Type password"

$tc002.Range("A9").Value = "Click"
$tc002.Range("B9").Value = "Xpath"
$tc002.Range("C9").Value = '//XCUIElementTypeSwitch[contains(@name, "Remember Login")]'
$tc002.Range("D9").Value = "N/A"
$tc002.Range("E9").Value = "Click in the Log in button"

$tc002.Range("A3").Select()

# add the hyperlink for D6 and restore its style index to the original Hyperlink style
$tc002.Hyperlinks.Add($tc002.Range("D6"), "mailto:test.admin@ur.com")
$tc002.Range("D6").Style = "Hyperlink"

# ---- 3. Sheet1 (Matrix) view tweaks ----
$matrix = $wb.Worksheets.Item("Matrix")
$matrix.Columns.Item(4).AutoFit()
$matrix.Range("D1").Select()
$matrix.Activate()

Write-Host "done"
